# Auto-generated edit script: Add data for 2023-03-06
# Applies per-cell numeric updates across Citywide Totals, By Neighborhood,
# and individual neighborhood sheets, per the diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 1118
$ws.Range("I3").Value = 7488
$ws.Range("J3").Value = 1199
$ws.Range("F4").Value = 1880
$ws.Range("I4").Value = 1752
$ws.Range("J4").Value = 263
$ws.Range("J5").Value = 90
$ws.Range("I6").Value = 8968
$ws.Range("J6").Value = 1596
$ws.Range("F7").Value = 24070
$ws.Range("I7").Value = 26191
$ws.Range("J7").Value = 4266

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J2").Value = 19
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 56

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 39
$ws.Range("J7").Value = 147

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J3").Value = 18
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 365
$ws.Range("J3").Value = 62
$ws.Range("J6").Value = 52
$ws.Range("I7").Value = 979
$ws.Range("J7").Value = 155

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I6").Value = 101
$ws.Range("I7").Value = 262

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 27
$ws.Range("J7").Value = 111

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J4").Value = 22
$ws.Range("J8").Value = 271
$ws.Range("J11").Value = 52
$ws.Range("J15").Value = 47
$ws.Range("J17").Value = 11
$ws.Range("J18").Value = 61
$ws.Range("J19").Value = 146
$ws.Range("J20").Value = 90
$ws.Range("J22").Value = 7
$ws.Range("J25").Value = 24
$ws.Range("J29").Value = 239
$ws.Range("I31").Value = 262
$ws.Range("J33").Value = 175
$ws.Range("J36").Value = 62
$ws.Range("J37").Value = 147
$ws.Range("J41").Value = 28
$ws.Range("J42").Value = 179
$ws.Range("J43").Value = 48
$ws.Range("J52").Value = 95
$ws.Range("J53").Value = 41
$ws.Range("J54").Value = 81
$ws.Range("J60").Value = 26
$ws.Range("F63").Value = 172
$ws.Range("I63").Value = 183
$ws.Range("J63").Value = 24
$ws.Range("J65").Value = 111
$ws.Range("I67").Value = 979
$ws.Range("J67").Value = 155
$ws.Range("J69").Value = 11
$ws.Range("J73").Value = 39
$ws.Range("J76").Value = 73
$ws.Range("J77").Value = 35
$ws.Range("J79").Value = 127
$ws.Range("J85").Value = 175
$ws.Range("J87").Value = 18
$ws.Range("J89").Value = 47
$ws.Range("J91").Value = 58
$ws.Range("J95").Value = 73
$ws.Range("J96").Value = 56
$ws.Range("J99").Value = 53
$ws.Range("F101").Value = 24070
$ws.Range("I101").Value = 26191
$ws.Range("J101").Value = 4266

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 73

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 40
$ws.Range("J3").Value = 48
$ws.Range("J7").Value = 175

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 13
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 81

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 67
$ws.Range("J6").Value = 63
$ws.Range("J7").Value = 239

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J6").Value = 58
$ws.Range("J7").Value = 146

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J2").Value = 10
$ws.Range("J3").Value = 17
$ws.Range("J7").Value = 73

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 44
$ws.Range("J3").Value = 63
$ws.Range("J6").Value = 49
$ws.Range("J7").Value = 175

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J3").Value = 6
$ws.Range("J7").Value = 28

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J6").Value = 107
$ws.Range("J7").Value = 179

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("J3").Value = 3
$ws.Range("J7").Value = 11

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J3").Value = 25
$ws.Range("J7").Value = 58

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 35
$ws.Range("J3").Value = 42
$ws.Range("J7").Value = 127

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J3").Value = 33
$ws.Range("J7").Value = 90

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 61

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("J2").Value = 2
$ws.Range("J7").Value = 11

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J2").Value = 19
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 95

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J6").Value = 5
$ws.Range("J7").Value = 24

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J5").Value = 3
$ws.Range("J7").Value = 52

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 39

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 89
$ws.Range("J3").Value = 88
$ws.Range("J7").Value = 271

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 26

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 5
$ws.Range("J6").Value = 24
$ws.Range("J7").Value = 41

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("J2").Value = 3
$ws.Range("J7").Value = 7

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J5").Value = 2
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 22

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 18
